$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.329.88'
$ws.Range("E2").Value = '  -0.22%  '
$ws.Range("D3").Value = '1.911.21'
$ws.Range("E3").Value = '  +0.13%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.718'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +8.83%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '254.43'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.09%  '
$ws.Range("E7").Value = '  +0.27%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '40.69'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.35%  '
$ws.Range("E9").Value = '  +2.15%  '
$ws.Range("E10").Value = '  -1.12%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0751'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.84%  '
$ws.Range("E12").Value = '  -0.65%  '
$ws.Range("D13").Value = '2.191.50'
$ws.Range("E13").Value = '  +0.44%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '12.63'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.59%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.719'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.00%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.92'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.00%  '
$ws.Range("D17").Value = '1.917.51'
$ws.Range("E17").Value = '  +0.66%  '
$ws.Range("D18").Value = '35.338.83'
$ws.Range("E18").Value = '  -0.13%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '74.42'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.41%  '
$ws.Range("D20").Value = '0.0₃0851'
$ws.Range("E20").Value = '  +3.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '244.20'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.54%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.99'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.75%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.07'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.06%  '
$ws.Range("E24").Value = '  +0.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.47'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.79%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.36'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.61%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '166.55'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.36%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.65'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.74'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.41%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.132'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.61%  '
$ws.Range("D31").Value = '4.128.65'
$ws.Range("E31").Value = '  +19.45%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.36'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.07%  '
$ws.Range("E33").Value = '  +15.09%  '
$ws.Range("E34").Value = '  +22.12%  '
$ws.Range("E35").Value = '  +3.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.21'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.92%  '
$ws.Range("E37").Value = '  +0.24%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.915'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.14%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.03'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.08%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0220'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.00%  '
$ws.Range("B41").Value = 'InjectiveProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '17.31'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.07%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '96.99'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +7.40%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.11'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.86%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0648'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.73%  '
$ws.Range("D45").Value = '1.339.89'
$ws.Range("E45").Value = '  -0.35%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.43'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.25%  '
$ws.Range("E47").Value = '  +0.64%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.77'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.97%  '
$ws.Range("E49").Value = '  -1.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '45.21'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.42%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '11.84'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +12.60%  '
